{"js": "// Apply red (\"EE0000\") font color to the \"13 dicembre\" day-header\n// paragraph's run AND paragraph mark, so it matches the formatting\n// already used by every other day-header paragraph in the document\n// (e.g. \"1 dicembre\", \"12 dicembre\", \"15 dicembre\", ...).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = \"13 dicembre\";\nlet found = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.trim() === target) {\n    found = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!found) {\n  throw new Error(`Paragraph \"${target}\" not found`);\n}\n\n// Setting Font.color on the paragraph (not just its range of runs)\n// colors both the run text and the paragraph mark, matching the other\n// headers' <w:pPr><w:rPr><w:color .../></w:rPr></w:pPr> + run <w:rPr>.\nfound.font.color = \"#EE0000\";\n\nawait context.sync();\n", "ps1": "# Apply red (\"EE0000\") font color to the \"13 dicembre\" day-header\n# paragraph's run AND paragraph mark, so it matches the formatting\n# already used by every other day-header paragraph in the document\n# (e.g. \"1 dicembre\", \"12 dicembre\", \"15 dicembre\", ...).\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.Trim() -eq \"13 dicembre\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Paragraph '13 dicembre' not found\"\n}\n\n# Setting Font.Color on the paragraph's Range (which includes the\n# paragraph mark) colors both the run text and the paragraph mark,\n# matching the <w:pPr><w:rPr><w:color .../></w:rPr></w:pPr> + run\n# <w:rPr><w:color .../></w:rPr> pattern used by the other headers.\n$target.Range.Font.Color = 238\n"}
